$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 616.5
$ws.Range("I4").Value = 539.8
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 539.8
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -425.8
$ws.Range("N4").Value = -1228

$ws.Range("H5").Value = 96.22221999999999
$ws.Range("I5").Value = 49
$ws.Range("J5").Value = 190.66667
$ws.Range("K5").Value = 49
$ws.Range("L5").Value = 190.66667
$ws.Range("M5").Value = 66
$ws.Range("N5").Value = -420.66667

$ws.Range("H6").Value = 337.5
$ws.Range("I6").Value = 383.33334
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 1150.00002
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -1038.00002
$ws.Range("N6").Value = -824

$ws.Range("H8").Value = 4444
$ws.Range("I8").Value = 4444
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 13332
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -13193
$ws.Range("N8").Value = ""

$ws.Range("H9").Value = 211.66667
$ws.Range("I9").Value = 50
$ws.Range("K9").Value = 50
$ws.Range("M9").Value = 119

$ws.Range("H10").Value = 50000
$ws.Range("I10").Value = 50000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 50000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -49707
$ws.Range("N10").Value = ""

$ws.Range("H12").Value = 167.125
$ws.Range("I12").Value = 172.83333
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 172.83333
$ws.Range("L12").Value = 150
$ws.Range("M12").Value = -2.833329999999989
$ws.Range("N12").Value = -490

$ws.Range("H18").Value = 340.5
$ws.Range("I18").Value = 340.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 340.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -56.5
$ws.Range("N18").Value = ""

$ws.Range("H107").Value = 483736.9
$ws.Range("I107").Value = 855216.9399999999
$ws.Range("J107").Value = 812.9
$ws.Range("K107").Value = 855216.9399999999
$ws.Range("L107").Value = 812.9
$ws.Range("M107").Value = -853296.9399999999
$ws.Range("N107").Value = -4652.9

$ws.Range("H132").Value = 224918.16
$ws.Range("I132").Value = 253821.06
$ws.Range("J132").Value = 51500.75
$ws.Range("K132").Value = 761463.1799999999
$ws.Range("L132").Value = 154502.25
$ws.Range("M132").Value = -758933.1799999999
$ws.Range("N132").Value = -159562.25

$ws.Range("H135").Value = 2491.6875
$ws.Range("I135").Value = 2206.6924
$ws.Range("K135").Value = 19860.2316
$ws.Range("M135").Value = -17325.2316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 25155
$ws.Range("I97").Value = 25155
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 25155
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -24659
$ws.Range("N97").Value = ""

$ws.Range("H139").Value = 48000
$ws.Range("J139").Value = 48000
$ws.Range("L139").Value = 48000
$ws.Range("N139").Value = -58280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 43000
$ws.Range("J59").Value = 43000
$ws.Range("L59").Value = 43000
$ws.Range("N59").Value = -44694

$ws.Range("H86").Value = 7023.1665
$ws.Range("I86").Value = 1510.8
$ws.Range("J86").Value = 13913.625
$ws.Range("K86").Value = 1510.8
$ws.Range("L86").Value = 13913.625
$ws.Range("M86").Value = -387.8
$ws.Range("N86").Value = -16159.625

$ws.Range("H89").Value = 7023.1665
$ws.Range("I89").Value = 1510.8
$ws.Range("J89").Value = 13913.625
$ws.Range("K89").Value = 7554
$ws.Range("L89").Value = 69568.125
$ws.Range("M89").Value = -1938
$ws.Range("N89").Value = -80800.125

$ws.Range("H105").Value = 2560.9
$ws.Range("I105").Value = 2401
$ws.Range("K105").Value = 2401
$ws.Range("M105").Value = -654

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1619.1364
$ws.Range("I31").Value = 1223.4445
$ws.Range("J31").Value = 3399.75
$ws.Range("K31").Value = 1223.4445
$ws.Range("L31").Value = 3399.75
$ws.Range("M31").Value = -928.4445000000001
$ws.Range("N31").Value = -3989.75

$ws.Range("H34").Value = 1619.1364
$ws.Range("I34").Value = 1223.4445
$ws.Range("J34").Value = 3399.75
$ws.Range("K34").Value = 1223.4445
$ws.Range("L34").Value = 3399.75
$ws.Range("M34").Value = -1021.4445
$ws.Range("N34").Value = -3803.75

$ws.Range("H58").Value = 1931.0769
$ws.Range("I58").Value = 1398.1666
$ws.Range("J58").Value = 2387.8572
$ws.Range("K58").Value = 1398.1666
$ws.Range("L58").Value = 2387.8572
$ws.Range("M58").Value = -1195.1666
$ws.Range("N58").Value = -2793.8572

$ws.Range("H136").Value = 1931.0769
$ws.Range("I136").Value = 1398.1666
$ws.Range("J136").Value = 2387.8572
$ws.Range("K136").Value = 4194.4998
$ws.Range("L136").Value = 7163.571599999999
$ws.Range("M136").Value = -1644.4998
$ws.Range("N136").Value = -12263.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 612.5333000000001
$ws.Range("I92").Value = 350
$ws.Range("J92").Value = 842.25
$ws.Range("K92").Value = 1050
$ws.Range("L92").Value = 2526.75
$ws.Range("M92").Value = 198
$ws.Range("N92").Value = -5022.75

$ws.Range("H113").Value = 11364372
$ws.Range("I113").Value = 540.6
$ws.Range("J113").Value = 12821274
$ws.Range("K113").Value = 1621.8
$ws.Range("L113").Value = 38463822
$ws.Range("M113").Value = 548.1999999999998
$ws.Range("N113").Value = -38468162

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 867.7778
$ws.Range("I97").Value = 734.6667
$ws.Range("J97").Value = 1533.3334
$ws.Range("K97").Value = 734.6667
$ws.Range("L97").Value = 1533.3334
$ws.Range("M97").Value = -238.6667
$ws.Range("N97").Value = -2525.3334

$ws.Range("H107").Value = 314.3
$ws.Range("I107").Value = 311.25
$ws.Range("K107").Value = 311.25
$ws.Range("M107").Value = 1608.75

$ws.Range("H126").Value = 2374.95
$ws.Range("I126").Value = 2166.3333
$ws.Range("J126").Value = 2411.7646
$ws.Range("K126").Value = 6498.999899999999
$ws.Range("L126").Value = 7235.293799999999
$ws.Range("M126").Value = -4028.999899999999
$ws.Range("N126").Value = -12175.2938

$ws.Range("H132").Value = 2499.5312
$ws.Range("I132").Value = 2079.7307
$ws.Range("J132").Value = 4318.6665
$ws.Range("K132").Value = 6239.1921
$ws.Range("L132").Value = 12955.9995
$ws.Range("M132").Value = -3709.1921
$ws.Range("N132").Value = -18015.9995

$ws.Range("H137").Value = 43500
$ws.Range("J137").Value = 43500
$ws.Range("L137").Value = 43500
$ws.Range("N137").Value = -53700

$ws.Range("H138").Value = 55800
$ws.Range("J138").Value = 55800
$ws.Range("L138").Value = 55800
$ws.Range("N138").Value = -66080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""
